# coordinadoresTorresRegiones.xlsx
# Commit: "agregado PabloSanMartin, actualizado gitignore"
#
# MARTIN PALMA (IdCoordinador YP11856) is replaced by
# PABLO SEBASTIAN SAN MARTIN (IdCoordinador SE45933) for the two CENTRONORTE
# rows (COMUNICACIONES and INFRAESTRUCTURA) in which he was the coordinador.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: CENTRONORTE / COMUNICACIONES
$ws.Range("C2").Value = "PABLO SEBASTIAN SAN MARTIN"
$ws.Range("D2").Value = "SE45933"

# Row 3: CENTRONORTE / INFRAESTRUCTURA
$ws.Range("C3").Value = "PABLO SEBASTIAN SAN MARTIN"
$ws.Range("D3").Value = "SE45933"

# Leave the same cell selected as in the authored workbook
$ws.Range("D3").Select()
